# Split the "Chai tea: The spice of life" tagline run into several runs
# separated by w:proofErr gramStart/gramEnd markers (the grammar checker
# flagging "spice" and "life"), matching what Word itself writes out when
# it re-saves a document after running the grammar checker over this
# sentence fragment. The paragraph's own properties (style/numbering) are
# left untouched -- only the run content inside the paragraph changes.

$d = $word.ActiveDocument

# Locate the target run's text and remember its extent.
$target = $d.Content
$found = $target.Find.Execute("Chai tea: The spice of life", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Chai tea: The spice of life' tagline"
}

# Insert the new run/proofErr sequence right after the existing run; Word
# keeps the paragraph's own identity/properties and simply appends the new
# runs supplied in the wrapper <w:p>.
$newRunsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Chai tea: The </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>spice</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> of </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>life</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'

$target.InsertXML($newRunsXml)

# Now remove the original (now redundant) run by re-finding its exact text
# and deleting it -- Delete() clears only the run content, leaving the
# paragraph mark, its paraId/rsid attributes and <w:pPr> untouched.
$old = $d.Content
$foundOld = $old.Find.Execute("Chai tea: The spice of life", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $foundOld) {
    throw "Could not find the original tagline run to remove"
}

$old.Delete()
